$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.673704333333333
$ws.Range("N2").Value = 29.021113
$ws.Range("O2").Value = 0.1714456165911166
$ws.Range("P2").Value = 0.1714456165911166
$ws.Range("Q2").Value = 100.9377849487211
$ws.Range("R2").Value = 908.4400645384899
$ws.Range("S2").Value = 0.166492415207697
$ws.Range("T2").Value = 0.166492415207697

# Row 3
$ws.Range("O3").Value = 0.4148961799842911
$ws.Range("P3").Value = 0.4148961799842911
$ws.Range("S3").Value = 0.4029094965476718
$ws.Range("T3").Value = 0.4029094965476718

# Row 4
$ws.Range("M4").Value = 23.09142233333333
$ws.Range("N4").Value = 69.27426699999999
$ws.Range("O4").Value = 0.4092458280188166
$ws.Range("P4").Value = 0.4092458280188166
$ws.Range("Q4").Value = 240.9415195387677
$ws.Range("R4").Value = 2168.47367584891
$ws.Range("S4").Value = 0.397422387782745
$ws.Range("T4").Value = 0.397422387782745

# Row 5
$ws.Range("M5").Value = 0.2489653333333333
$ws.Range("N5").Value = 0.746896
$ws.Range("O5").Value = 0.004412375405775742
$ws.Range("P5").Value = 0.004412375405775742
$ws.Range("Q5").Value = 2.597764869564445
$ws.Range("R5").Value = 23.37988382608
$ws.Range("S5").Value = 0.004284898341044606
$ws.Range("T5").Value = 0.004284898341044606

# Row 6
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.673704333333333
$ws.Range("N6").Value = 29.021113
$ws.Range("O6").Value = 0.1714456165911166
$ws.Range("P6").Value = 0.1714456165911166
$ws.Range("Q6").Value = 3.002930646562
$ws.Range("R6").Value = 27.026375819058
$ws.Range("S6").Value = 0.004953201383419629
$ws.Range("T6").Value = 0.00495320138341963

# Row 7
$ws.Range("O7").Value = 0.4148961799842911
$ws.Range("P7").Value = 0.4148961799842911
$ws.Range("S7").Value = 0.01198668343661924
$ws.Range("T7").Value = 0.01198668343661924

# Row 8
$ws.Range("M8").Value = 23.09142233333333
$ws.Range("N8").Value = 69.27426699999999
$ws.Range("O8").Value = 0.4092458280188166
$ws.Range("P8").Value = 0.4092458280188166
$ws.Range("Q8").Value = 7.168085503557998
$ws.Range("R8").Value = 64.51276953202199
$ws.Range("S8").Value = 0.01182344023607161
$ws.Range("T8").Value = 0.01182344023607161

# Row 9
$ws.Range("M9").Value = 0.2489653333333333
$ws.Range("N9").Value = 0.746896
$ws.Range("O9").Value = 0.004412375405775742
$ws.Range("P9").Value = 0.004412375405775742
$ws.Range("Q9").Value = 0.077284316704
$ws.Range("R9").Value = 0.695558850336
$ws.Range("S9").Value = 0.0001274770647311351
$ws.Range("T9").Value = 0.0001274770647311352
